# Git Basic Commands slide (slide 5) - tidy up the trailing commas/spaces in
# the bulleted list of common git commands inside the Subtitle placeholder.
#
# Before (runs, each its own paragraph bullet):
#   "git clone, "   -> "git clone"
#   "git status, "  -> "git status "
#   "git add,"      -> "git add"
#   " git commit, " -> " git commit "
#   "git push, "    -> "git push"
#
# We target each run's text precisely via TextRange.Characters(start, length)
# so unrelated runs/paragraphs are left untouched. Edits are applied from the
# end of the text backwards so earlier character offsets stay valid.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(5)
$shape = $s.Shapes.Item(2)
$tr = $shape.TextFrame.TextRange

# "git push, " (start=76, length=10) -> "git push"
$tr.Characters(76, 10).Text = "git push"

# " git commit, " (start=62, length=13) -> " git commit "
$tr.Characters(62, 13).Text = " git commit "

# "git add," (start=53, length=8) -> "git add"
$tr.Characters(53, 8).Text = "git add"

# "git status, " (start=40, length=12) -> "git status "
$tr.Characters(40, 12).Text = "git status "

# "git clone, " (start=28, length=11) -> "git clone"
$tr.Characters(28, 11).Text = "git clone"
